$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to text so strings like "1.010" or "26.966.16"
# are stored verbatim instead of being auto-coerced to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '26.966.16'
$ws.Range('E2').Value = '  +1.90%  '

$ws.Range('D3').Value = '1.818.10'
$ws.Range('E3').Value = '  +2.52%  '

$ws.Range('D4').Value = '1.010'
$ws.Range('E4').Value = '  +0.69%  '

$ws.Range('D5').Value = '312.41'
$ws.Range('E5').Value = '  +1.94%  '

$ws.Range('D6').Value = '1.007'
$ws.Range('E6').Value = '  +0.54%  '

$ws.Range('D7').Value = '0.4301'
$ws.Range('E7').Value = '  +0.31%  '

$ws.Range('D8').Value = '0.3677'
$ws.Range('E8').Value = '  +0.40%  '

$ws.Range('D9').Value = '0.07248'
$ws.Range('E9').Value = '  +0.17%  '

$ws.Range('D10').Value = '2.132.33'
$ws.Range('E10').Value = '  +18.41%  '

$ws.Range('D11').Value = '0.8668'
$ws.Range('E11').Value = '  +2.21%  '

$ws.Range('D12').Value = '21.31'
$ws.Range('E12').Value = '  +4.77%  '

$ws.Range('D13').Value = '5.420'
$ws.Range('E13').Value = '  +2.99%  '

$ws.Range('D14').Value = '6.618'
$ws.Range('E14').Value = '  +2.83%  '

$ws.Range('D15').Value = '0.06984'
$ws.Range('E15').Value = '  +2.37%  '

$ws.Range('D16').Value = '81.20'
$ws.Range('E16').Value = '  +2.05%  '

$ws.Range('D17').Value = '1.013'
$ws.Range('E17').Value = '  +0.79%  '

$ws.Range('D18').Value = '0.000008895'
$ws.Range('E18').Value = '  +2.28%  '

$ws.Range('D19').Value = '1.007'
$ws.Range('E19').Value = '  +0.44%  '

$ws.Range('D20').Value = '15.22'
$ws.Range('E20').Value = '  +1.07%  '

$ws.Range('D21').Value = '27.023.57'
$ws.Range('E21').Value = '  +2.12%  '

$ws.Range('D22').Value = '5.202'
$ws.Range('E22').Value = '  +1.95%  '

$ws.Range('D23').Value = '11.00'
$ws.Range('E23').Value = '  -2.36%  '

$ws.Range('D24').Value = '2.356.72'
$ws.Range('E24').Value = '  +16.76%  '

$ws.Range('D25').Value = '154.17'
$ws.Range('E25').Value = '  +0.97%  '

$ws.Range('D26').Value = '1.890'
$ws.Range('E26').Value = '  +2.20%  '

$ws.Range('E27').Value = '  +1.33%  '

$ws.Range('D28').Value = '5.226'
$ws.Range('E28').Value = '  +2.35%  '

$ws.Range('D29').Value = '1.906'
$ws.Range('E29').Value = '  +11.64%  '

$ws.Range('D30').Value = '114.89'
$ws.Range('E30').Value = '  +0.28%  '

$ws.Range('D31').Value = '0.08961'
$ws.Range('E31').Value = '  +0.14%  '

$ws.Range('D32').Value = '1.186'
$ws.Range('E32').Value = '  +5.95%  '

$ws.Range('D33').Value = '0.7493'
$ws.Range('E33').Value = '  +2.85%  '

$ws.Range('D34').Value = '4.429'
$ws.Range('E34').Value = '  +1.85%  '

$ws.Range('D35').Value = '2.814'
$ws.Range('E35').Value = '  +2.19%  '

$ws.Range('D36').Value = '1.006'
$ws.Range('E36').Value = '  +0.37%  '

$ws.Range('D37').Value = '1.133'
$ws.Range('E37').Value = '  +4.95%  '

$ws.Range('D38').Value = '0.05234'
$ws.Range('E38').Value = '  +1.57%  '

$ws.Range('D39').Value = '0.01925'
$ws.Range('E39').Value = '  +1.48%  '

$ws.Range('D40').Value = '0.5109'
$ws.Range('E40').Value = '  +3.59%  '

$ws.Range('D41').Value = '0.1656'
$ws.Range('E41').Value = '  +2.60%  '

$ws.Range('D42').Value = '2.747'

$ws.Range('D43').Value = '6.473'
$ws.Range('E43').Value = '  +3.90%  '

$ws.Range('D44').Value = '8.349'

$ws.Range('D45').Value = '106.95'
$ws.Range('E45').Value = '  +1.92%  '

$ws.Range('D46').Value = '10.43'
$ws.Range('E46').Value = '  +2.61%  '

$ws.Range('D47').Value = '1.007'
$ws.Range('E47').Value = '  +0.56%  '

$ws.Range('D48').Value = '0.4589'
$ws.Range('E48').Value = '  +1.88%  '

$ws.Range('D49').Value = '1.648'
$ws.Range('E49').Value = '  +3.93%  '

$ws.Range('D50').Value = '0.06226'
$ws.Range('E50').Value = '  +0.35%  '

$ws.Range('D51').Value = '1.830'
$ws.Range('E51').Value = '  +4.57%  '

# Restore default (Normal/General) styling on column D so no visible
# formatting change is introduced beyond the text content itself.
$ws.Range("D2:D51").Style = "Normal"
